$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.919399999999992
$ws.Range("B7").Value = 5.7591
$ws.Range("A10").Value = -22.1423
$ws.Range("E10").Value = 16.48139999999999
$ws.Range("A12").Value = -21.59370000000001
$ws.Range("E14").Value = 16.7944
$ws.Range("B15").Value = 4.605899999999995
$ws.Range("A18").Value = -22.40520000000002
$ws.Range("D18").Value = -8.16109999999999
$ws.Range("D19").Value = -8.956499999999991
$ws.Range("B20").Value = 9.799999999999992
$ws.Range("D27").Value = -9.169700000000004
$ws.Range("B29").Value = 4.891500000000003
$ws.Range("B30").Value = 4.687500000000003
$ws.Range("B31").Value = 4.966100000000001
$ws.Range("E32").Value = 16.55319999999999
$ws.Range("E35").Value = 16.361
$ws.Range("A37").Value = -19.984
$ws.Range("B40").Value = 9.098099999999995
$ws.Range("D42").Value = -8.718699999999998
$ws.Range("E43").Value = 17.24
$ws.Range("D44").Value = -7.761399999999997
$ws.Range("D47").Value = -7.5483
$ws.Range("E49").Value = 15.48219999999999
$ws.Range("A55").Value = -22.424
$ws.Range("E56").Value = 16.61160000000001
$ws.Range("D58").Value = -8.27249999999999
$ws.Range("A68").Value = -21.5466
$ws.Range("B68").Value = 4.712300000000002
$ws.Range("E69").Value = 17.39860000000002
$ws.Range("D73").Value = -7.911699999999994
$ws.Range("B76").Value = 5.846799999999997
$ws.Range("A77").Value = -20.59979999999999
$ws.Range("A78").Value = -20.01529999999998
$ws.Range("E81").Value = 16.52219999999999
$ws.Range("B87").Value = 4.990199999999999
$ws.Range("B88").Value = 4.695099999999999
$ws.Range("E92").Value = 18.33980000000002
$ws.Range("D95").Value = -7.994
$ws.Range("B96").Value = 4.896200000000004
$ws.Range("B98").Value = 5.841299999999998
$ws.Range("B101").Value = 9.174199999999997
$ws.Range("D101").Value = -7.743100000000001
$ws.Range("B102").Value = 8.696000000000005
